# Commit: "modelling CO2 emissions and relevant emission tax"
#
# Adds a new "Emission_CO2_sink" / "emission_CO2" node__commodity relationship
# row to the two rel_for_node_basic_structure(_ptdf) sheets, and a matching
# "emission_CO2" commodity_physics_none row to obj_commodity_ptdf.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: rel_for_node_basic_structure
# Insert a new data row right after the header (row 2), i.e. at row 3,
# pushing the existing data rows down by one.
# ---------------------------------------------------------------------------
$wsStruct = $wb.Worksheets.Item("rel_for_node_basic_structure")
$wsStruct.Rows.Item(3).Insert()
$wsStruct.Range("A3").Value = "node__commodity"
$wsStruct.Range("B3").Value = "node__stochastic_structure"
$wsStruct.Range("C3").Value = "node__temporal_block"
$wsStruct.Range("D3").Value = "Emission_CO2_sink"
$wsStruct.Range("E3").Value = "emission_CO2"
$wsStruct.Range("F3").Value = "default"
$wsStruct.Range("G3").Value = "blk_t1"
$wsStruct.Range("H3").Value = "blk_t2"
$wsStruct.Range("A3:H3").Select()

# ---------------------------------------------------------------------------
# Sheet: rel_for_node_basic_struc_ptdf
# Row 3 here is the special "ptdf_group" row (no A3/E3 values); clear its A3
# value, then insert the new Emission_CO2_sink row right after it (row 4).
# ---------------------------------------------------------------------------
$wsPtdf = $wb.Worksheets.Item("rel_for_node_basic_struc_ptdf")
$wsPtdf.Range("A3").ClearContents()
$wsPtdf.Rows.Item(4).Insert()
$wsPtdf.Range("A4").Value = "node__commodity"
$wsPtdf.Range("B4").Value = "node__stochastic_structure"
$wsPtdf.Range("C4").Value = "node__temporal_block"
$wsPtdf.Range("D4").Value = "Emission_CO2_sink"
$wsPtdf.Range("E4").Value = "emission_CO2"
$wsPtdf.Range("F4").Value = "default"
$wsPtdf.Range("G4").Value = "blk_t1"
$wsPtdf.Range("H4").Value = "blk_t2"

# ---------------------------------------------------------------------------
# Sheet: obj_commodity_ptdf
# Append a new "emission_CO2" commodity row.
# ---------------------------------------------------------------------------
$wsCommodity = $wb.Worksheets.Item("obj_commodity_ptdf")
$wsCommodity.Range("A4").Value = "commodity"
$wsCommodity.Range("B4").Value = "emission_CO2"
$wsCommodity.Range("C4").Value = "commodity_physics"
$wsCommodity.Range("D4").Value = "commodity_physics_none"
$wsCommodity.Activate()
$wsCommodity.Range("G4").Select()

# ---------------------------------------------------------------------------
# Restore the originally active sheet/selection: rel_for_node_basic_struc_ptdf
# with cell A3 selected (now the empty "ptdf_group" A cell).
# ---------------------------------------------------------------------------
$wsPtdf.Activate()
$wsPtdf.Range("A3").Select()
